# "new updates for split drive"
# Clear the leftover formula-text values in C3:D3 (replace with numeric 0)
# and correct D4 from 1 back to 0, then move the active selection to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("C4").Select()
